# Fruta / hortaliza, semanal
#
# A new weekly price record was added for "Vega Modelo de Temuco - Poroto
# granado". It is inserted as a new row 43 (pushing the previously-existing
# rows 43-59 down to 44-60), keeping the data sorted the same way it was
# before (newest entry first in this block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 43, shifting rows 43:59 down to 44:60.
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new record.
$ws.Cells.Item(43, 1).Value  = 10
$ws.Cells.Item(43, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(43, 3).Value  = "La Araucanía"
$ws.Cells.Item(43, 4).Value  = 44582
$ws.Cells.Item(43, 5).Value  = 9
$ws.Cells.Item(43, 6).Value  = 100112030
$ws.Cells.Item(43, 7).Value  = "Poroto granado"
$ws.Cells.Item(43, 8).Value  = "Sin especificar"
$ws.Cells.Item(43, 9).Value  = "Primera"
$ws.Cells.Item(43, 10).Value = 50
$ws.Cells.Item(43, 11).Value = 28000
$ws.Cells.Item(43, 12).Value = 28000
$ws.Cells.Item(43, 13).Value = 28000
$ws.Cells.Item(43, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(43, 15).Value = "Región del Maule"
$ws.Cells.Item(43, 16).Value = 1120
$ws.Cells.Item(43, 17).Value = 25
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# Match the date-number-format style used by the other rows in column D.
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
